$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cluster names in cross-join order
$clusters = @("ECs", "FAPs", "sCs")
$ligand = "Alcam"
$receptor = "Chl1"

# Numeric data for each (sending, target) cluster pair, in row order
# Columns: E F G H I J K L M N O P Q R S T
$data = @{
    "ECs|ECs"   = @(3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 2, 0.6666666666666666, 0.07850066666666666, 0.235502, 0.01562049898812527, 0.01562049898812527, 3.303184728786, 29.728662559074, 0.01527229555455865, 0.01527229555455865)
    "ECs|FAPs"  = @(3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 2, 0.6666666666666666, 0.051033, 0.153099, 0.010154829999673, 0.010154829999673, 2.147388467157, 19.326496204413, 0.009928464204581593, 0.009928464204581594)
    "ECs|sCs"   = @(3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 3, 1, 4.895956666666667, 14.68787, 0.9742246710122017, 0.9742246710122018, 206.01416498541, 1854.12748486869, 0.9525077991139579, 0.952507799113958)
    "FAPs|ECs"  = @(2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 2, 0.6666666666666666, 0.07850066666666666, 0.235502, 0.01562049898812527, 0.01562049898812527, 0.03745735193977778, 0.337116167458, 0.0001731843043866494, 0.0001731843043866495)
    "FAPs|FAPs" = @(2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 2, 0.6666666666666666, 0.051033, 0.153099, 0.010154829999673, 0.010154829999673, 0.024350889269, 0.219158003421, 0.0001125864910586392, 0.0001125864910586392)
    "FAPs|sCs"  = @(2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 3, 1, 4.895956666666667, 14.68787, 0.9742246710122017, 0.9742246710122018, 2.336153051081111, 21.02537745973, 0.01080121845619798, 0.01080121845619798)
    "sCs|ECs"   = @(3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 2, 0.6666666666666666, 0.07850066666666666, 0.235502, 0.01562049898812527, 0.01562049898812527, 0.03785419897666666, 0.34068779079, 0.0001750191291799772, 0.0001750191291799772)
    "sCs|FAPs"  = @(3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 2, 0.6666666666666666, 0.051033, 0.153099, 0.010154829999673, 0.010154829999673, 0.024608878095, 0.221479902855, 0.0001137793040327697, 0.0001137793040327697)
    "sCs|sCs"   = @(3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 3, 1, 4.895956666666667, 14.68787, 0.9742246710122017, 0.9742246710122018, 2.360903744016667, 21.24813369615, 0.01091565344204598, 0.01091565344204598)
}

# Expand the 3x3 cross join of sending/target clusters into rows 2..10,
# replacing the 2 old data rows and adding 7 new ones.
$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sending
        $ws.Cells.Item($row, 2).Value = $ligand
        $ws.Cells.Item($row, 3).Value = $receptor
        $ws.Cells.Item($row, 4).Value = $target

        $vals = $data["$sending|$target"]
        for ($i = 0; $i -lt $vals.Length; $i++) {
            $ws.Cells.Item($row, 5 + $i).Value = $vals[$i]
        }

        $row++
    }
}
